$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the erroneous numeric values from the table (rows 3-12, columns C-L)
# replacing them with a single-space placeholder string, matching the
# blank-row convention used further down the table (rows 13-21).
$rows = 3..12
$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L")

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = " "
    }
}

# D10 uses the double-space placeholder instead of the single-space one
$ws.Range("D10").Value = "  "

# Update the active selection to match the saved view state
$ws.Range("E7").Select()
